$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New phone-number link text replacing the old one shown in A3
$ws.Range("A3").Value = "https://api.whatsapp.com/send?phone=919542856170"

# B3 gets the same weather report message that already lives in B2
$ws.Cells.Item(3, 2).Value = "CURRENT WEATHER: 29°C`nRealFeel® 35°`nTONIGHT’S WEATHER FORECAST: 29°Lo`nRealFeel® 33°`nPlenty of clouds"
$ws.Cells.Item(3, 2).Style = "Normal"

# Column B needs to be widened now that it holds the full phone-number / message text
$ws.Columns("B").ColumnWidth = 35.5

$ws.Range("B3").Select() | Out-Null
